$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update NumEp (column D) values
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("D22").Value = 5
$ws.Range("D25").Value = 2
$ws.Range("D26").Value = 3

# Widen column E (target stored width 11.42578125 chars; nearest value the
# engine's pixel-quantized ColumnWidth setter can reach is 11.5)
$ws.Columns.Item(5).ColumnWidth = 10.666666666666668

# Update the view: scroll back to A1 (remove frozen/top-left offset) and change selection to D6
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D6").Select()
